$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsInteger    = $wb.Worksheets.Item("Integer")
$wsBoolean    = $wb.Worksheets.Item("Boolean")
$wsSubscript  = $wb.Worksheets.Item("Subscript")

# ---------------------------------------------------------------------------
# "Boolean" sheet: the BVTQaZ and VTQaZ input-file rows each get split out
# into per-vehicle-type files (LDVs / HDVs / aircraft / rail / ships /
# motorbikes) instead of a single combined csv.
# ---------------------------------------------------------------------------

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" -> becomes the first of
# six rows.
$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
for ($k = 0; $k -lt 5; $k++) {
    $wsBoolean.Rows.Item(18).Insert()
}
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the insert above, "trans/VTQaZ/VTQaZ.csv" (originally row 21) is now
# at row 26 -> becomes the first of six rows too.
$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
for ($k = 0; $k -lt 5; $k++) {
    $wsBoolean.Rows.Item(27).Insert()
}
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of formatted-but-empty trailing rows (33-38) were left below the
# data on this sheet; touch-and-clear each one so the used range grows to
# match without leaving stray values behind.
for ($r = 33; $r -le 38; $r++) {
    $wsBoolean.Cells.Item($r, 1).Value = "x"
    $wsBoolean.Cells.Item($r, 1).ClearContents()
}

# ---------------------------------------------------------------------------
# View/selection bookkeeping to match the saved state: "Boolean" was scrolled
# down with A32 selected, "Integer" had A13 selected, and "About" ended up
# the active tab.
# ---------------------------------------------------------------------------
$wsInteger.Range("A13").Select()
$wsBoolean.Range("A32").Select()
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
